# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the zh-cn and de-de report sheets (regenerated report timestamps).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-16 09:56:13"
$wsZhCn.Range("G2").Value = "2016-02-16 09:57:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-16 09:56:26"
$wsDeDe.Range("G2").Value = "2016-02-16 09:57:32"
